# "all food and meal data from external JSON"
# The "Food Items" lookup table was refreshed from the external JSON source:
# a few items were removed (Greek Yoghurt plain, Brown Sugar 1 tsp, Jam,
# Full-Fat Milk, Peas) and new ones were added (Nespresso coffee pod, Polish
# Wafer), while the rest kept their values. The "Meals" sheet itself still
# references the same foods/codes/counts/totals as before - those numbers
# are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Food Items")

$foods = @(
    @("Frozen Strawberries", "FSXX", 100, "g", 32),
    @("Frozen Bananas", "FBXX", 100, "g", 89),
    @("Lidl Greek Style Yoghurt (Full Fat)", "LGSY", 100, "g", 126),
    @("Brown Sugar (1 tbsp)", "BS1T", 1, "tbsp", 52),
    @("Semi-skimmed Milk (50 ml)", "SM5M", 50, "ml", 25),
    @("Honey", "HONE", 1, "tbsp", 64),
    @("Homemade Vegetable Soup", "HVSX", 1, "bowl", 100),
    @("Sourdough Bread", "SBXX", 1, "slice", 174),
    @("Butter", "BUTT", 1, "tbsp", 102),
    @("Nespresso coffee pod", "NCPX", 1, "pod", 1),
    @("Semi-skimmed Milk (100 ml)", "SM1M", 100, "ml", 50),
    @("Polish Wafer", "PWXX", 1, "wafer", 49)
)

for ($i = 0; $i -lt $foods.Length; $i++) {
    $row = $i + 2
    $item = $foods[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}

# The old table ran through row 16; the refreshed table only needs 12 rows
# (through row 13), so remove what's left over.
$ws.Range("A14:E16").Delete()
